$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# Copy the date/time formatting of the last existing data row (41) down
# onto the three new rows being appended (42:44), matching columns A & B.
$ws.Range("A41:B41").Copy()
$ws.Range("A42:B44").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# New data rows appended to the log (rows 42-44)
$ws.Range("A42").Value = 43697
$ws.Range("B42").Value = 0.57152777777777775
$ws.Range("C42").Value = 72
$ws.Range("D42").Value = 21888

$ws.Range("A43").Value = 43697
$ws.Range("B43").Value = 0.63194444444444442
$ws.Range("C43").Value = 72
$ws.Range("D43").Value = 21920
$ws.Range("E43").Value = 17413

$ws.Range("A44").Value = 43697
$ws.Range("B44").Value = 0.67361111111111116
$ws.Range("C44").Value = 72
$ws.Range("D44").Value = 22261

# Update view state to match the scrolled/selected position after the edit
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("C44").Select()
